# PpGoGreen.pptx - very small changes:
# Prefix the titles of the "Stats", "Achievements" and "Personal Settings"
# slides with "Application: ".

$p = $ppt.ActivePresentation

# Slide 10 ("Stats") - add a new, separate run "Application: " before the
# existing "Stats" run.
$s10 = $p.Slides.Item(10)
$title10 = $s10.Shapes.Item(1).TextFrame.TextRange
$title10.InsertBefore("Application: ")

# Slide 11 ("Achievements") - add a new, separate run "Application: " before
# the existing "Achievements" run.
$s11 = $p.Slides.Item(11)
$title11 = $s11.Shapes.Item(1).TextFrame.TextRange
$title11.InsertBefore("Application: ")

# Slide 12 ("Personal Settings") - prepend "Application: " inside the
# existing "Personal " run (so it stays a single run reading
# "Application: Personal ").
$s12 = $p.Slides.Item(12)
$title12 = $s12.Shapes.Item(1).TextFrame.TextRange
$personalRun = $title12.Characters(1, 9)
$personalRun.InsertBefore("Application: ")
